$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textRefs = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($r in $textRefs) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range("D2").Value = "55.234.84"
$ws.Range("E2").Value = "  -5.22%  "
$ws.Range("D3").Value = "2.881.48"
$ws.Range("E3").Value = "  -6.04%  "
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "484.40"
$ws.Range("E5").Value = "  -7.36%  "
$ws.Range("D6").Value = "131.07"
$ws.Range("E6").Value = "  -8.15%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.414"
$ws.Range("E8").Value = "  -7.44%  "
$ws.Range("D9").Value = "7.08"
$ws.Range("E9").Value = "  -5.98%  "
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  -8.71%  "
$ws.Range("D11").Value = "0.342"
$ws.Range("E11").Value = "  -7.69%  "
$ws.Range("D12").Value = "3.354.85"
$ws.Range("E12").Value = "  -5.36%  "
$ws.Range("D13").Value = "0.124"
$ws.Range("E13").Value = "  -5.05%  "
$ws.Range("D14").Value = "25.41"
$ws.Range("E14").Value = "  -5.64%  "
$ws.Range("D15").Value = "0.0000155"
$ws.Range("E15").Value = "  -9.65%  "
$ws.Range("D16").Value = "55.074.78"
$ws.Range("E16").Value = "  -5.52%  "
$ws.Range("D17").Value = "5.88"
$ws.Range("E17").Value = "  -5.80%  "
$ws.Range("D18").Value = "2.864.44"
$ws.Range("E18").Value = "  -6.46%  "
$ws.Range("D19").Value = "12.23"
$ws.Range("E19").Value = "  -6.55%  "
$ws.Range("D20").Value = "7.55"
$ws.Range("E20").Value = "  -7.42%  "
$ws.Range("D21").Value = "309.73"
$ws.Range("E21").Value = "  -8.48%  "
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "5.77"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "0.472"
$ws.Range("E24").Value = "  -6.25%  "
$ws.Range("D25").Value = "61.14"
$ws.Range("E25").Value = "  -6.48%  "
$ws.Range("D26").Value = "0.980"
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("D27").Value = "0.158"
$ws.Range("E27").Value = "  -6.37%  "
$ws.Range("D28").Value = "0.0₃0820"
$ws.Range("E28").Value = "  -14.75%  "
$ws.Range("D29").Value = "6.27"
$ws.Range("E29").Value = "  -9.83%  "
$ws.Range("D30").Value = "6.88"
$ws.Range("E30").Value = "  -9.40%  "
$ws.Range("D31").Value = "1.71"
$ws.Range("E31").Value = "  -7.24%  "
$ws.Range("D32").Value = "19.45"
$ws.Range("E32").Value = "  -8.09%  "
$ws.Range("D33").Value = "1.10"
$ws.Range("E33").Value = "  -9.91%  "
$ws.Range("D34").Value = "146.49"
$ws.Range("E34").Value = "  -7.14%  "
$ws.Range("D35").Value = "4.32"
$ws.Range("E35").Value = "  -10.15%  "
$ws.Range("D36").Value = "5.52"
$ws.Range("E36").Value = "  -7.77%  "
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").Value = "23.85"
$ws.Range("E37").Value = "  -6.20%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.17"
$ws.Range("E38").Value = "  -10.15%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.0644"
$ws.Range("E39").Value = "  -7.49%  "
$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").Value = "2.907.28"
$ws.Range("E40").Value = "  -6.36%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "0.995"
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("D42").Value = "35.71"
$ws.Range("E42").Value = "  -5.45%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.59"
$ws.Range("E43").Value = "  -8.36%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "0.623"
$ws.Range("E44").Value = "  -6.63%  "
$ws.Range("D45").Value = "2.076.29"
$ws.Range("E45").Value = "  -11.23%  "
$ws.Range("D46").Value = "1.30"
$ws.Range("E46").Value = "  -10.79%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "0.914"
$ws.Range("E47").Value = "  -10.51%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "5.75"
$ws.Range("E48").Value = "  -5.54%  "
$ws.Range("D49").Value = "0.0224"
$ws.Range("E49").Value = "  -7.43%  "
$ws.Range("D50").Value = "18.39"
$ws.Range("E50").Value = "  -7.47%  "
$ws.Range("D51").Value = "0.0828"
$ws.Range("E51").Value = "  -8.30%  "
